$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# stay as text (matching the source inline-string cells) by temporarily
# switching to a text number format, then restoring the default style so
# no stray formatting is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.585.95"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "1.584.12"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "210.94"
$ws.Range("E5").Value = "  -2.48%  "
Set-TextValue "D6" "0.506"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  -0.11%  "
Set-TextValue "D8" "0.249"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("E9").Value = "  -1.04%  "
Set-TextValue "D10" "19.52"
$ws.Range("E10").Value = "  -3.88%  "
Set-TextValue "D11" "0.0832"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Value = "1.805.61"
$ws.Range("D13").Value = "1.583.64"
$ws.Range("E13").Value = "  -2.94%  "
Set-TextValue "D14" "4.04"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  -2.47%  "
Set-TextValue "D16" "64.45"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "26.609.58"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("E18").Value = "  -0.51%  "
Set-TextValue "D19" "208.42"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("E20").Value = "  -0.07%  "
Set-TextValue "D21" "6.74"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("E24").Value = "  -2.26%  "
Set-TextValue "D25" "146.35"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -4.18%  "
Set-TextValue "D30" "0.0501"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  -3.74%  "
Set-TextValue "D33" "0.670"
$ws.Range("E33").Value = "  +23.87%  "
Set-TextValue "D34" "2.95"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").Value = "1.319.81"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("E38").Value = "  -1.26%  "
Set-TextValue "D39" "0.823"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.32"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.18"
$ws.Range("E43").Value = "  -3.68%  "
Set-TextValue "D44" "63.27"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "1.719.19"
$ws.Range("E45").Value = "  -2.83%  "
Set-TextValue "D46" "89.07"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("E47").Value = "  +1.00%  "
Set-TextValue "D48" "0.837"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -0.30%  "
